$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.426.11'
$ws.Range('E2').Value = '  +4.97%  '
$ws.Range('D3').Value = '3.504.06'
$ws.Range('E3').Value = '  +2.23%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'418.06"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('D6').Value = "'131.77"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.41%  '
$ws.Range('D7').Value = "'0.655"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.45%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  +6.50%  '
$ws.Range('D10').Value = "'0.162"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +15.36%  '
$ws.Range('D11').Value = "'43.11"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').Value = "'0.0000266"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +20.90%  '
$ws.Range('D13').Value = "'9.95"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.41%  '
$ws.Range('D14').Value = '4.062.98'
$ws.Range('E14').Value = '  +2.37%  '
$ws.Range('D15').Value = "'0.140"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').Value = "'20.40"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '3.500.99'
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').Value = "'12.79"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('D20').Value = '65.403.93'
$ws.Range('E20').Value = '  +5.02%  '
$ws.Range('D21').Value = "'449.05"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.42%  '
$ws.Range('D22').Value = "'90.04"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('E23').Value = '  -0.79%  '
$ws.Range('D24').Value = "'13.23"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('E25').Value = '  +2.93%  '
$ws.Range('D26').Value = "'9.88"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.47%  '
$ws.Range('D27').Value = "'33.94"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.31%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = "'12.48"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.95%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'2.72"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.23%  '
$ws.Range('D30').Value = "'7.41"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.00%  '
$ws.Range('D31').Value = "'0.117"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.46%  '
$ws.Range('E32').Value = '  -1.57%  '
$ws.Range('D33').Value = "'0.999"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').Value = "'39.37"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.34%  '
$ws.Range('D35').Value = "'57.28"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('E36').Value = '  +4.00%  '
$ws.Range('D37').Value = '0.0₃0733'
$ws.Range('E37').Value = '  +35.37%  '
$ws.Range('D38').Value = "'0.147"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.84%  '
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('E40').Value = '  +0.95%  '
$ws.Range('D41').Value = "'2.77"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.25%  '
$ws.Range('D42').Value = "'4.49"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.41%  '
$ws.Range('D43').Value = "'145.97"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.49%  '
$ws.Range('D44').Value = "'3.27"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('E45').Value = '  -3.48%  '
$ws.Range('D46').Value = "'1.98"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.59%  '
$ws.Range('E47').Value = '  -2.46%  '
$ws.Range('D48').Value = "'15.74"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.35%  '
$ws.Range('E49').Value = '  +3.56%  '
$ws.Range('E50').Value = '  +10.31%  '
$ws.Range('D51').Value = "'21.59"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.08%  '
